# Applies the "Updated cryptos list" data refresh (GitHub Actions scrape)
# to Sheet1. Each data row (2-51) holds a coin's Name (col B), Link (col C),
# Price (col D) and Volume(1h) (col E) as plain text. This script rewrites
# only the cells whose values actually changed between runs, including the
# two pairs of rows (43/44 and 50/51) whose coins were swapped/reordered
# along with their Name, Link, Price and Volume values.
#
# Price values (col D) are plain numeric-looking strings (e.g. "512.30")
# that Excel would otherwise silently coerce to numbers (and round, e.g.
# 512.30 -> 512.29999999999995) when assigned through Range.Value. To keep
# them as genuine text - matching how the source data is stored - each D
# cell is briefly switched to Text format ("@") before the assignment and
# then restored to the workbook's default "Normal" style so no visible
# formatting/style change is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("D2")
$rng.NumberFormat = "@"
$rng.Value = "56.911.65"
$rng.Style = "Normal"
$ws.Range("E2").Value = "  -3.68%  "
$rng = $ws.Range("D3")
$rng.NumberFormat = "@"
$rng.Value = "2.536.31"
$rng.Style = "Normal"
$ws.Range("E3").Value = "  -4.76%  "
$ws.Range("E4").Value = "  -0.01%  "
$rng = $ws.Range("D5")
$rng.NumberFormat = "@"
$rng.Value = "512.30"
$rng.Style = "Normal"
$ws.Range("E5").Value = "  -2.48%  "
$rng = $ws.Range("D6")
$rng.NumberFormat = "@"
$rng.Value = "139.27"
$rng.Style = "Normal"
$ws.Range("E6").Value = "  -3.61%  "
$ws.Range("E7").Value = "  +0.03%  "
$rng = $ws.Range("D8")
$rng.NumberFormat = "@"
$rng.Value = "0.553"
$rng.Style = "Normal"
$ws.Range("E8").Value = "  -2.86%  "
$ws.Range("E9").Value = "  -7.50%  "
$rng = $ws.Range("D10")
$rng.NumberFormat = "@"
$rng.Value = "0.0987"
$rng.Style = "Normal"
$ws.Range("E10").Value = "  -4.18%  "
$rng = $ws.Range("D11")
$rng.NumberFormat = "@"
$rng.Value = "0.322"
$rng.Style = "Normal"
$ws.Range("E11").Value = "  -4.08%  "
$ws.Range("E12").Value = "  -0.27%  "
$rng = $ws.Range("D13")
$rng.NumberFormat = "@"
$rng.Value = "2.981.11"
$rng.Style = "Normal"
$ws.Range("E13").Value = "  -4.85%  "
$rng = $ws.Range("D14")
$rng.NumberFormat = "@"
$rng.Value = "56.934.16"
$rng.Style = "Normal"
$ws.Range("E14").Value = "  -3.63%  "
$rng = $ws.Range("D15")
$rng.NumberFormat = "@"
$rng.Value = "19.97"
$rng.Style = "Normal"
$ws.Range("E15").Value = "  -5.44%  "
$ws.Range("E16").Value = "  -3.38%  "
$rng = $ws.Range("D17")
$rng.NumberFormat = "@"
$rng.Value = "2.538.20"
$rng.Style = "Normal"
$ws.Range("E17").Value = "  -4.98%  "
$rng = $ws.Range("D18")
$rng.NumberFormat = "@"
$rng.Value = "330.89"
$rng.Style = "Normal"
$ws.Range("E18").Value = "  -2.40%  "
$ws.Range("E19").Value = "  -2.90%  "
$rng = $ws.Range("D20")
$rng.NumberFormat = "@"
$rng.Value = "10.04"
$rng.Style = "Normal"
$ws.Range("E20").Value = "  -3.57%  "
$rng = $ws.Range("D21")
$rng.NumberFormat = "@"
$rng.Value = "6.11"
$rng.Style = "Normal"
$ws.Range("E21").Value = "  -4.54%  "
$ws.Range("E22").Value = "  +0.14%  "
$rng = $ws.Range("D23")
$rng.NumberFormat = "@"
$rng.Value = "64.10"
$rng.Style = "Normal"
$ws.Range("E23").Value = "  -0.58%  "
$rng = $ws.Range("D24")
$rng.NumberFormat = "@"
$rng.Value = "0.164"
$rng.Style = "Normal"
$ws.Range("E24").Value = "  -1.10%  "
$rng = $ws.Range("D25")
$rng.NumberFormat = "@"
$rng.Value = "0.998"
$rng.Style = "Normal"
$ws.Range("E25").Value = "  +0.05%  "
$rng = $ws.Range("D26")
$rng.NumberFormat = "@"
$rng.Value = "0.398"
$rng.Style = "Normal"
$ws.Range("E26").Value = "  -5.00%  "
$rng = $ws.Range("D27")
$rng.NumberFormat = "@"
$rng.Value = "2.654.77"
$rng.Style = "Normal"
$ws.Range("E27").Value = "  -4.49%  "
$rng = $ws.Range("D28")
$rng.NumberFormat = "@"
$rng.Value = "6.86"
$rng.Style = "Normal"
$ws.Range("E28").Value = "  -3.31%  "
$rng = $ws.Range("D29")
$rng.NumberFormat = "@"
$rng.Value = "0.0₃0746"
$rng.Style = "Normal"
$ws.Range("E29").Value = "  -7.02%  "
$ws.Range("E30").Value = "  -0.04%  "
$rng = $ws.Range("D31")
$rng.NumberFormat = "@"
$rng.Value = "6.23"
$rng.Style = "Normal"
$ws.Range("E31").Value = "  -6.86%  "
$ws.Range("E32").Value = "  -3.03%  "
$rng = $ws.Range("D33")
$rng.NumberFormat = "@"
$rng.Value = "147.98"
$rng.Style = "Normal"
$ws.Range("E33").Value = "  -1.96%  "
$rng = $ws.Range("D34")
$rng.NumberFormat = "@"
$rng.Value = "18.39"
$rng.Style = "Normal"
$ws.Range("E34").Value = "  -2.61%  "
$rng = $ws.Range("D35")
$rng.NumberFormat = "@"
$rng.Value = "3.95"
$rng.Style = "Normal"
$ws.Range("E35").Value = "  -4.97%  "
$ws.Range("E36").Value = "  -5.92%  "
$ws.Range("E37").Value = "  -5.70%  "
$rng = $ws.Range("D38")
$rng.NumberFormat = "@"
$rng.Value = "35.66"
$rng.Style = "Normal"
$ws.Range("E38").Value = "  -3.39%  "
$ws.Range("E39").Value = "  -6.22%  "
$ws.Range("E40").Value = "  -3.10%  "
$ws.Range("E41").Value = "  +0.02%  "
$rng = $ws.Range("D42")
$rng.NumberFormat = "@"
$rng.Value = "3.45"
$rng.Style = "Normal"
$ws.Range("E42").Value = "  -3.86%  "
$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$rng = $ws.Range("D43")
$rng.NumberFormat = "@"
$rng.Value = "10.61"
$rng.Style = "Normal"
$ws.Range("E43").Value = "  -0.47%  "
$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$rng = $ws.Range("D44")
$rng.NumberFormat = "@"
$rng.Value = "0.0949"
$rng.Style = "Normal"
$ws.Range("E44").Value = "  -2.23%  "
$rng = $ws.Range("D45")
$rng.NumberFormat = "@"
$rng.Value = "0.574"
$rng.Style = "Normal"
$ws.Range("E45").Value = "  -6.93%  "
$rng = $ws.Range("D46")
$rng.NumberFormat = "@"
$rng.Value = "258.53"
$rng.Style = "Normal"
$ws.Range("E46").Value = "  -6.20%  "
$rng = $ws.Range("D47")
$rng.NumberFormat = "@"
$rng.Value = "0.0516"
$rng.Style = "Normal"
$ws.Range("E47").Value = "  -2.96%  "
$rng = $ws.Range("D48")
$rng.NumberFormat = "@"
$rng.Value = "18.39"
$rng.Style = "Normal"
$ws.Range("E48").Value = "  -7.62%  "
$rng = $ws.Range("D49")
$rng.NumberFormat = "@"
$rng.Value = "1.963.19"
$rng.Style = "Normal"
$ws.Range("E49").Value = "  -3.98%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$rng = $ws.Range("D50")
$rng.NumberFormat = "@"
$rng.Value = "0.0219"
$rng.Style = "Normal"
$ws.Range("E50").Value = "  -4.34%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$rng = $ws.Range("D51")
$rng.NumberFormat = "@"
$rng.Value = "4.50"
$rng.Style = "Normal"
$ws.Range("E51").Value = "  -4.50%  "
